$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected NN values for rows 2-9 (columns B..H)
# Row 2
$ws.Range("B2").Value = 1202.253173828125
$ws.Range("C2").Value = 0.9019
$ws.Range("D2").Value = 0.8992999792098999
$ws.Range("E2").Value = 1.041900038719177
$ws.Range("F2").Value = 0.7164000272750854
$ws.Range("H2").Value = 0.6542

# Row 3
$ws.Range("B3").Value = 1213.76318359375
$ws.Range("C3").Value = 0.971
$ws.Range("D3").Value = 0.9056999999999999
$ws.Range("E3").Value = 1.574699997901917
$ws.Range("F3").Value = 0.822700023651123
$ws.Range("H3").Value = 0.7107

# Row 4
$ws.Range("B4").Value = 842.8361206054688
$ws.Range("C4").Value = 0.9846
$ws.Range("D4").Value = 0.9523
$ws.Range("E4").Value = 1.442999958992004
$ws.Range("F4").Value = 0.8457000255584717
$ws.Range("H4").Value = 1.1233

# Row 5
$ws.Range("B5").Value = 908.736328125
$ws.Range("C5").Value = 0.9586
$ws.Range("D5").Value = 0.9213
$ws.Range("E5").Value = 1.223899960517883
$ws.Range("F5").Value = 0.8320000171661377
$ws.Range("H5").Value = 0.8486

# Row 6
$ws.Range("B6").Value = 1183.89404296875
$ws.Range("C6").Value = 0.9366
$ws.Range("D6").Value = 0.918
$ws.Range("E6").Value = 1.178400039672852
$ws.Range("F6").Value = 0.8572999835014343
$ws.Range("H6").Value = 0.8195

# Row 7
$ws.Range("B7").Value = 922.9119873046875
$ws.Range("C7").Value = 0.9294
$ws.Range("D7").Value = 0.9118000268936157
$ws.Range("E7").Value = 1.16129994392395
$ws.Range("F7").Value = 0.8374999761581421
$ws.Range("H7").Value = 0.7648

# Row 8
$ws.Range("B8").Value = 1019.223388671875
$ws.Range("C8").Value = 0.9133
$ws.Range("D8").Value = 0.9073
$ws.Range("E8").Value = 1.029000043869019
$ws.Range("F8").Value = 0.8392000198364258
$ws.Range("H8").Value = 0.7248

# Row 9
$ws.Range("B9").Value = 7293.6181640625
$ws.Range("C9").Value = 0.9399
$ws.Range("D9").Value = 0.9121
$ws.Range("E9").Value = 1.574699997901917
$ws.Range("F9").Value = 0.7164000272750854
$ws.Range("H9").Value = 5.6459
